$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 870
$ws1.Range("F5").Value = 77

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 870
$ws4.Range("F6").Value = 77
